# Automated data refresh: update extraction timestamps and latest readings
# pulled from meteo.cat for 2026-02-21 (run at 19:20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell (never edited by this run) used as a formatting donor so that
# percentage values re-entered as literal text ("NN%") keep cell style s=3
# instead of Excel re-interpreting them as a numeric percentage.
$fmtDonor = $ws.Range("H2")

$ws.Range("E2").Value = "2026-02-21 19:18:31"
$ws.Range("E3").Value = "2026-02-21 19:18:33"
$ws.Range("O3").Value = "1.9 °C"
$ws.Range("E4").Value = "2026-02-21 19:18:36"
$ws.Range("E5").Value = "2026-02-21 19:18:38"
$ws.Range("H5").Value = "'39%"
$fmtDonor.Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("O5").Value = "3.6 °C"
$ws.Range("E6").Value = "2026-02-21 19:18:41"
$ws.Range("H6").Value = "'69%"
$fmtDonor.Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("E7").Value = "2026-02-21 19:18:43"
$ws.Range("E8").Value = "2026-02-21 19:18:46"
$ws.Range("E9").Value = "2026-02-21 19:18:48"
$ws.Range("H9").Value = "'53%"
$fmtDonor.Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("N9").Value = "8.5 °C 18:59 TU"
$ws.Range("O9").Value = "13.9 °C"
$ws.Range("E10").Value = "2026-02-21 19:18:51"
$ws.Range("H10").Value = "'78%"
$fmtDonor.Copy()
$ws.Range("H10").PasteSpecial(-4122)
$ws.Range("E11").Value = "2026-02-21 19:18:53"
$ws.Range("H11").Value = "'50%"
$fmtDonor.Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("O11").Value = "9.3 °C"
$ws.Range("E12").Value = "2026-02-21 19:18:56"
$ws.Range("H12").Value = "'59%"
$fmtDonor.Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("O12").Value = "13.0 °C"
$ws.Range("E13").Value = "2026-02-21 19:18:58"
$ws.Range("O13").Value = "5.4 °C"
$ws.Range("E14").Value = "2026-02-21 19:19:01"
$ws.Range("H14").Value = "'67%"
$fmtDonor.Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("O14").Value = "11.8 °C"
$ws.Range("E15").Value = "2026-02-21 19:19:03"
$ws.Range("H15").Value = "'52%"
$fmtDonor.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("O15").Value = "13.7 °C"
$ws.Range("E16").Value = "2026-02-21 19:19:06"
$ws.Range("E17").Value = "2026-02-21 19:19:08"
$ws.Range("H17").Value = "'33%"
$fmtDonor.Copy()
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("E18").Value = "2026-02-21 19:19:10"
$ws.Range("E19").Value = "2026-02-21 19:19:13"
$ws.Range("H19").Value = "'64%"
$fmtDonor.Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("E20").Value = "2026-02-21 19:19:15"
$ws.Range("E21").Value = "2026-02-21 19:19:18"
$ws.Range("E22").Value = "2026-02-21 19:19:20"
$ws.Range("E23").Value = "2026-02-21 19:19:22"
$ws.Range("O23").Value = "2.8 °C"
$ws.Range("E24").Value = "2026-02-21 19:19:25"
$ws.Range("K24").Value = "15.8 MJ/m2"
$ws.Range("E25").Value = "2026-02-21 19:19:28"
$ws.Range("E26").Value = "2026-02-21 19:19:30"
$ws.Range("E27").Value = "2026-02-21 19:19:33"
$ws.Range("H27").Value = "'33%"
$fmtDonor.Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("O27").Value = "4.3 °C"
$ws.Range("E28").Value = "2026-02-21 19:19:35"
$ws.Range("E29").Value = "2026-02-21 19:19:37"
$ws.Range("H29").Value = "'64%"
$fmtDonor.Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("K29").Value = "15.3 MJ/m2"
$ws.Range("O29").Value = "12.0 °C"
$ws.Range("E30").Value = "2026-02-21 19:19:40"
$ws.Range("H30").Value = "'64%"
$fmtDonor.Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("O30").Value = "11.9 °C"
$ws.Range("E31").Value = "2026-02-21 19:19:42"
$ws.Range("E32").Value = "2026-02-21 19:19:45"
$ws.Range("H32").Value = "'78%"
$fmtDonor.Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("O32").Value = "5.7 °C"
$ws.Range("E33").Value = "2026-02-21 19:19:47"
$ws.Range("H33").Value = "'53%"
$fmtDonor.Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("J33").Value = "1030.3 hPa"
$ws.Range("E34").Value = "2026-02-21 19:19:50"
$ws.Range("N34").Value = "0.3 °C 18:46 TU"
$ws.Range("O34").Value = "4.9 °C"
$ws.Range("E35").Value = "2026-02-21 19:19:53"
$ws.Range("J35").Value = "1030.8 hPa"
$ws.Range("E36").Value = "2026-02-21 19:19:55"
$ws.Range("O36").Value = "13.6 °C"
$ws.Range("E37").Value = "2026-02-21 19:19:57"
$ws.Range("E38").Value = "2026-02-21 19:20:00"
$ws.Range("E39").Value = "2026-02-21 19:20:03"
$ws.Range("E40").Value = "2026-02-21 19:20:05"
$ws.Range("J40").Value = "1030.4 hPa"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-21 19:20:07"
$ws.Range("O41").Value = "11.4 °C"
$ws.Range("E42").Value = "2026-02-21 19:20:10"
$ws.Range("H42").Value = "'72%"
$fmtDonor.Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("E43").Value = "2026-02-21 19:20:12"
$ws.Range("H43").Value = "'75%"
$fmtDonor.Copy()
$ws.Range("H43").PasteSpecial(-4122)
$ws.Range("O43").Value = "7.1 °C"
$ws.Range("E44").Value = "2026-02-21 19:20:15"
$ws.Range("O44").Value = "2.5 °C"
$ws.Range("E45").Value = "2026-02-21 19:20:17"
$ws.Range("E46").Value = "2026-02-21 19:20:20"
$ws.Range("J46").Value = "1031.5 hPa"

$excel.CutCopyMode = $false
